# Commit: "minor updates to illumina barcoding protocol"
#
# The underlying XML diff only touches the internal <w:nsid w:val="..."/>
# identifiers carried on four list definitions (w:abstractNum) inside
# word/numbering.xml:
#
#   abstractNumId 990   : 90cddb9b -> 7a703ff8
#   abstractNumId 99411 : b2564cfd -> 7883db44
#   abstractNumId 991   : 5fd9dd06 -> 9992cf74
#   abstractNumId 99414 : 7a36daca -> 3d0e1c90
#
# Nothing else in those list definitions (multilevel type, level
# formatting, numbering id mappings, paragraph content, etc.) changes.
# `nsid` is just Word's internal bookkeeping id for a list definition;
# it isn't surfaced anywhere in the Word object model (no
# ListTemplate/List/AbstractNum property reads or writes it - it is not
# even round-trippable through Document.WordOpenXML, which is read-only).
# We still perform the edit the way a real automation script would if it
# had raw-XML access, so that if the host exposes a writable surface for
# it, the swap is applied; the attempts are wrapped defensively so a
# locked-down host just leaves the document untouched instead of failing
# the whole script.

$d = $word.ActiveDocument

$nsidMap = @{
    "90cddb9b" = "7a703ff8";
    "b2564cfd" = "7883db44";
    "5fd9dd06" = "9992cf74";
    "7a36daca" = "3d0e1c90"
}

function Swap-Nsids([string]$xml) {
    $result = $xml
    foreach ($old in $nsidMap.Keys) {
        $new = $nsidMap[$old]
        $result = $result.Replace("<w:nsid w:val=`"$old`"/>", "<w:nsid w:val=`"$new`"/>")
        $result = $result.Replace("<w:nsid w:val=`"$old`" />", "<w:nsid w:val=`"$new`" />")
    }
    return $result
}

$applied = $false

# Preferred route: round-trip the package-level OOXML (Flat OPC) that
# Document.WordOpenXML exposes, patch just the four nsid values, and
# write it back. On hosts where this property is read-only this throws
# and we fall through to the other attempts below.
try {
    $full = $d.WordOpenXML
    if ($full -and $full.Contains("<w:nsid")) {
        $patched = Swap-Nsids $full
        if ($patched -ne $full) {
            $d.WordOpenXML = $patched
            $applied = $true
        }
    }
} catch {
    $applied = $false
}

if (-not $applied) {
    try {
        $full = $d.Content.WordOpenXML
        if ($full -and $full.Contains("<w:nsid")) {
            $patched = Swap-Nsids $full
            if ($patched -ne $full) {
                $d.Content.WordOpenXML = $patched
                $applied = $true
            }
        }
    } catch {
        $applied = $false
    }
}

if (-not $applied) {
    try {
        $full = $d.Content.XML
        if ($full -and $full.Contains("<w:nsid")) {
            $patched = Swap-Nsids $full
            if ($patched -ne $full) {
                $d.Content.XML = $patched
                $applied = $true
            }
        }
    } catch {
        $applied = $false
    }
}

if (-not $applied) {
    try {
        $word.Selection.WholeStory()
        $full = $word.Selection.WordOpenXML
        if ($full -and $full.Contains("<w:nsid")) {
            $patched = Swap-Nsids $full
            if ($patched -ne $full) {
                $word.Selection.WordOpenXML = $patched
                $applied = $true
            }
        }
    } catch {
        $applied = $false
    }
}

if ($applied) {
    Write-Output "nsid values patched via raw OOXML round-trip"
} else {
    Write-Output "nsid values are internal-only (not exposed by the Word object model); no reachable list formatting/content changed, so the document is left as-is"
}
